$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "und die Umstände" -> "und Umstände" (occurs once in the TOC
#    entry and once in the actual heading text). Find/Replace covers
#    both occurrences while keeping the existing (bold) run
#    formatting of the matched text.
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "und die Umstände", $false, $false, $false, $false, $false,
    $true, 1, $false, "und Umstände", 2)

# ------------------------------------------------------------------
# 2) "Einkäufen," -> "Einkäufe," inside the "Dabei ist das
#    Transportieren ..." sentence.
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "wie zum Beispiel Einkäufen,", $false, $false, $false, $false, $false,
    $true, 1, $false, "wie zum Beispiel Einkäufe,", 2)

# ------------------------------------------------------------------
# 3) Insert a brand-new paragraph right after the "... aufbringen
#    kann." paragraph, containing a sentence about online services
#    not matching the needs of older people.
# ------------------------------------------------------------------
$found = $d.Content
$null = $found.Find.Execute(
    "das man selten aufbringen kann.", $false, $false, $false, $false,
    $false, $true, 1, $false, "", 0)
$insertionPoint = $found.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

$newRange = $insertionPoint.Paragraphs(1).Next().Range
$newRange.Collapse(1)
$newRange.InsertBefore("Auch moderne Online-Dienstleistungen entsprechen häufig nicht den Bedürfnissen älterer Menschen, da diese sich in den wenigsten Fällen ausreichend mit Computern auskennen.")

# ------------------------------------------------------------------
# 4) Replace the whole "Unsere Lösung dafür besteht darin, ..." body
#    paragraph with the new text about loneliness driving the
#    companion-robot idea, then split it into two paragraphs: one
#    about social contact, one about the robot's transport function.
# ------------------------------------------------------------------
$oldRng = $d.Content
$null = $oldRng.Find.Execute(
    "Unsere Lösung dafür besteht darin, einen autonom fahrenden Roboter zu entwickeln, welcher eine die Einkäufe für einen transportiert, sodass man selbst Nichts schweres mehr tragen muss.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $oldRng.Duplicate
$target.Text = "Da viele ältere Menschen einsam leben, möchten sie meist nicht auf sozialen, zwischenmenschlichen Kontakt verzichten. Dieser Aspekt war für uns ein weiterer Grund, einen Roboter als „Gefährten“ des Besitzers einzusetzen. Der Roboter sollte unserer Vorstellung nach mithilfe einer Transportfläche beispielsweise Einkäufe tragen können."

# split that paragraph into two right after "... einzusetzen. "
$splitRng = $d.Content
$null = $splitRng.Find.Execute(
    "einzusetzen. Der Roboter sollte", $false, $false, $false, $false,
    $false, $true, 1, $false, "", 0)
$splitPoint = $splitRng.Duplicate
$splitPoint.Start = $splitPoint.Start + 13
$splitPoint.Collapse(1)
$splitPoint.InsertParagraphAfter()
